$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- Widen column C (Subjects) ---
$ws.Columns.Item(3).ColumnWidth = 50.125

# --- Add "Real" (actual finish) dates for rows 3-5 (col G) ---
$ws.Range("G3").Value = 41333
$ws.Range("G4").Value = 41332
$ws.Range("G5").Value = 41336

# --- Fill in row 9 (task #8: read/write file) ---
# Comments (J) first so the new shared strings land in the same order
# they appear in the target workbook (J9, then C9, then D9).
$ws.Range("J9").Value = "_Tạo đối tượng StreamWriter strW = new StreamWriter(@`"C:\test.txt`");`n_Gọi hàm của đối tượng vừa tạo strW.Write(chuỗi cần ghi viết ở đây)`n_Nhớ khai báo thư viện System.IO"
$ws.Range("J9").WrapText = $true

$ws.Range("C9").Value = "Tạo một form gồm 1 nút bấm button và một textbox. `nKhi bấm vào button thì chữ ở ô text sẽ được lưu xuống file có tên là C:\test.txt chẳng hạn"
$ws.Range("C9").WrapText = $true

$ws.Range("D9").Value = "_Làm quen với đối tượng ghi dữ liệu`n"
$ws.Range("D9").WrapText = $true

$ws.Range("E9").Value = 41397
$ws.Range("F9").Value = 41338

$ws.Range("H9").Value = "Pending"
$ws.Range("I9").Value = "Medium"

$ws.Rows.Item(9).RowHeight = 71.25

# --- Update the active selection to D6 (no more frozen top-left scroll) ---
$ws.Range("D6").Select()
